$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A42 contains a date-like text value ("2025-04-07"). Force it to be stored
# as text (matching the other date cells in the sheet) instead of letting
# Excel auto-convert it to a date serial number.
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "2025-04-07"
$ws.Range("A42").Style = "Normal"

$ws.Range("B42").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C42").Value = "NA"
$ws.Range("D42").Value = 1
